$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.517.42"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.394.27"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.15"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.35"
$ws.Range("E6").Value = "  -0.45%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.71"
$ws.Range("E9").Value = "  +2.04%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  -1.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.973.67"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.35"
$ws.Range("E14").Value = "  +0.23%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.399.89"
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.473.97"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.15"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.68"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "391.57"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "75.43"
$ws.Range("E22").Value = "  +1.45%  "
$ws.Range("E23").Value = "  -0.36%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000114"
$ws.Range("E25").Value = "  -3.37%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.193"
$ws.Range("E26").Value = "  +8.09%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.29"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.42"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.95"
$ws.Range("E34").Value = "  -0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.74"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("E36").Value = "  +1.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.426.54"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("E38").Value = "  -0.31%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0770"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "26.22"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.780"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.44"
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.15"
$ws.Range("E44").Value = "  +0.85%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.470.81"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "23.00"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.08"
$ws.Range("E50").Value = "  -0.76%  "
$ws.Range("E51").Value = "  -1.23%  "
